# Apply updated crypto price (D) and volume change (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.655.43"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").Value = "'1.889.44"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'236.97"
$ws.Range("E5").Value = "  -3.62%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "'0.4879"
$ws.Range("E7").Value = "  -2.68%  "
$ws.Range("D8").Value = "'0.2923"
$ws.Range("E8").Value = "  -1.82%  "
$ws.Range("D9").Value = "'0.06675"
$ws.Range("E9").Value = "  -2.26%  "
$ws.Range("D10").Value = "'1.889.75"
$ws.Range("E10").Value = "  -0.51%  "
$ws.Range("D11").Value = "'16.70"
$ws.Range("E11").Value = "  -2.43%  "
$ws.Range("D12").Value = "'0.07247"
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("D13").Value = "'89.17"
$ws.Range("E13").Value = "  -2.79%  "
$ws.Range("D14").Value = "'5.007"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("D15").Value = "'0.6636"
$ws.Range("E15").Value = "  -2.14%  "
$ws.Range("D16").Value = "'30.597.69"
$ws.Range("E16").Value = "  -0.84%  "
$ws.Range("D17").Value = "'0.000007905"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "'13.00"
$ws.Range("E19").Value = "  -1.88%  "
$ws.Range("D20").Value = "'2.134.91"
$ws.Range("E20").Value = "  -0.74%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").Value = "'4.748"
$ws.Range("E22").Value = "  -2.56%  "
$ws.Range("D23").Value = "'191.53"
$ws.Range("E23").Value = "  +5.10%  "
$ws.Range("D24").Value = "'6.096"
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("D25").Value = "'9.305"
$ws.Range("E25").Value = "  -0.49%  "
$ws.Range("D26").Value = "'159.89"
$ws.Range("E26").Value = "  +3.25%  "
$ws.Range("D27").Value = "'18.32"
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("D28").Value = "'1.828"
$ws.Range("E28").Value = "  -5.96%  "
$ws.Range("E29").Value = "  +0.69%  "
$ws.Range("D30").Value = "'4.259"
$ws.Range("E30").Value = "  -1.89%  "
$ws.Range("D31").Value = "'0.09019"
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("D32").Value = "'3.938"
$ws.Range("E32").Value = "  -3.00%  "
$ws.Range("E33").Value = "  -1.49%  "
$ws.Range("D34").Value = "'0.7320"
$ws.Range("E34").Value = "  -1.92%  "
$ws.Range("D35").Value = "'1.084"
$ws.Range("E35").Value = "  -4.52%  "
$ws.Range("D36").Value = "'2.688"
$ws.Range("E36").Value = "  +0.72%  "
$ws.Range("D38").Value = "'2.662"
$ws.Range("E38").Value = "  -2.32%  "
$ws.Range("D39").Value = "'0.9245"
$ws.Range("E39").Value = "  -1.36%  "
$ws.Range("D40").Value = "'2.044"
$ws.Range("E40").Value = "  -6.40%  "
$ws.Range("D41").Value = "'0.4406"
$ws.Range("E41").Value = "  +0.30%  "
$ws.Range("D42").Value = "'104.44"
$ws.Range("E42").Value = "  -1.38%  "
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").Value = "'5.735"
$ws.Range("E44").Value = "  -1.67%  "
$ws.Range("E45").Value = "  -0.76%  "
$ws.Range("D46").Value = "'7.347"
$ws.Range("E46").Value = "  -4.96%  "
$ws.Range("D47").Value = "'0.4122"
$ws.Range("E47").Value = "  +5.19%  "
$ws.Range("D48").Value = "'0.05826"
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("D49").Value = "'8.672"
$ws.Range("E49").Value = "  +1.10%  "
$ws.Range("D50").Value = "'1.414"
$ws.Range("E50").Value = "  +2.17%  "
$ws.Range("D51").Value = "'33.25"
$ws.Range("E51").Value = "  -0.20%  "
